$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Puerro" at Vega Modelo de Temuco.
# It belongs right after the existing row 111 (chronologically it is the
# most recent entry), so insert a fresh row at 112 - this pushes the old
# rows 112-135 down to 113-136, matching the dimension change to A1:R136.
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new record's data.
$ws.Range("A112").Value = 10
$ws.Range("B112").Value = "Vega Modelo de Temuco"
$ws.Range("C112").Value = "La Araucanía"
$ws.Range("D112").Value = 44476
$ws.Range("E112").Value = 9
$ws.Range("F112").Value = 100112005
$ws.Range("G112").Value = "Puerro"
$ws.Range("H112").Value = "Azul de Maquehue"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 50
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = 6000
$ws.Range("N112").Value = "$/docena de paquetes"
$ws.Range("O112").Value = "Provincia de Cautín"
$ws.Range("P112").Value = 500
$ws.Range("Q112").Value = 12
$ws.Range("R112").Value = "Hortaliza"
